$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Audio" / "audio.mp3" column (column A). This shifts
# columns B:G left to A:F, removing the old merged cell A2:A4 along
# with it (the merge was only on the deleted column).
$colA = $ws.Range("A1:A1048576")
$colA.Select()
$colA.EntireColumn.Delete()
